# Auto-generated edit script for cryptos.xlsx
# Applies updated Price (D) / Volume(1h) (E) figures, and restores the
# TrustWalletToken / TheSandbox row ordering that the source feed swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text even when it looks like a number
# (e.g. "327.99" or "1.003"), so it keeps matching the sheet's existing
# text-formatted Price/Coin/Link columns instead of becoming a numeric cell.
function Set-TextValue {
    param($Address, $Text)
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

$ws.Range('D2').Value = '27.475.74'
$ws.Range('E2').Value = '  -1.70%  '
$ws.Range('D3').Value = '1.746.22'
$ws.Range('E3').Value = '  -1.80%  '
$ws.Range('E4').Value = '  +0.11%  '
Set-TextValue 'D5' '327.99'
$ws.Range('E5').Value = '  +0.03%  '
Set-TextValue 'D6' '1.003'
$ws.Range('E6').Value = '  +0.00%  '
Set-TextValue 'D7' '0.4754'
$ws.Range('E7').Value = '  +5.48%  '
Set-TextValue 'D8' '0.3522'
$ws.Range('E8').Value = '  -1.77%  '
Set-TextValue 'D9' '42.69'
$ws.Range('E9').Value = '  +1.05%  '
Set-TextValue 'D10' '0.07454'
$ws.Range('E10').Value = '  -0.59%  '
Set-TextValue 'D11' '1.073'
$ws.Range('E11').Value = '  -2.69%  '
Set-TextValue 'D12' '1.003'
$ws.Range('E12').Value = '  +0.05%  '
Set-TextValue 'D13' '20.27'
$ws.Range('E13').Value = '  -3.45%  '
Set-TextValue 'D14' '6.022'
$ws.Range('E14').Value = '  -0.85%  '
$ws.Range('D15').Value = '1.744.70'
$ws.Range('E15').Value = '  -1.98%  '
Set-TextValue 'D16' '7.022'
$ws.Range('E16').Value = '  -2.98%  '
Set-TextValue 'D17' '91.83'
$ws.Range('E17').Value = '  -1.76%  '
Set-TextValue 'D18' '0.00001071'
$ws.Range('E18').Value = '  +0.57%  '
Set-TextValue 'D19' '0.06373'
$ws.Range('E19').Value = '  -1.43%  '
$ws.Range('E20').Value = '  +0.05%  '
Set-TextValue 'D21' '16.67'
$ws.Range('E21').Value = '  -3.39%  '
Set-TextValue 'D22' '5.749'
$ws.Range('E22').Value = '  -1.62%  '
$ws.Range('D23').Value = '27.552.63'
$ws.Range('E23').Value = '  -1.48%  '
Set-TextValue 'D24' '10.97'
$ws.Range('E24').Value = '  -3.64%  '
Set-TextValue 'D25' '2.153'
$ws.Range('E25').Value = '  +1.90%  '
Set-TextValue 'D26' '161.42'
$ws.Range('E26').Value = '  -1.29%  '
Set-TextValue 'D27' '19.83'
$ws.Range('E27').Value = '  -2.47%  '
$ws.Range('D28').Value = '1.947.73'
$ws.Range('E28').Value = '  -2.00%  '
Set-TextValue 'D29' '2.187'
$ws.Range('E29').Value = '  -1.92%  '
Set-TextValue 'D30' '121.32'
$ws.Range('E30').Value = '  -3.55%  '
Set-TextValue 'D31' '1.046'
$ws.Range('E31').Value = '  -5.66%  '
Set-TextValue 'D32' '0.09317'
$ws.Range('E32').Value = '  +1.62%  '
Set-TextValue 'D33' '3.621'
$ws.Range('E33').Value = '  -1.10%  '
Set-TextValue 'D34' '5.460'
$ws.Range('E34').Value = '  -2.34%  '
Set-TextValue 'D35' '0.02239'
$ws.Range('E35').Value = '  -2.87%  '
Set-TextValue 'D36' '11.31'
$ws.Range('E36').Value = '  -5.47%  '
Set-TextValue 'D37' '0.05930'
$ws.Range('E37').Value = '  -3.55%  '
Set-TextValue 'D38' '0.2043'
$ws.Range('E38').Value = '  -2.71%  '
$ws.Range('E39').Value = '  -3.33%  '
Set-TextValue 'D40' '1.434'
$ws.Range('E40').Value = '  +2.19%  '
$ws.Range('B41').Value = 'TheSandbox'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D41' '0.6064'
$ws.Range('E41').Value = '  -4.83%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue 'D42' '1.173'
$ws.Range('E42').Value = '  -1.62%  '
Set-TextValue 'D43' '7.698'
$ws.Range('E43').Value = '  -3.87%  '
Set-TextValue 'D44' '3.724'
$ws.Range('E44').Value = '  -0.59%  '
Set-TextValue 'D45' '12.90'
$ws.Range('E45').Value = '  -3.11%  '
Set-TextValue 'D46' '0.5689'
$ws.Range('E46').Value = '  -4.16%  '
Set-TextValue 'D47' '122.48'
$ws.Range('E47').Value = '  -0.66%  '
Set-TextValue 'D48' '1.901'
$ws.Range('E48').Value = '  -3.33%  '
Set-TextValue 'D49' '1.131'
$ws.Range('E49').Value = '  -1.28%  '
Set-TextValue 'D50' '0.06748'
$ws.Range('E50').Value = '  -2.67%  '
Set-TextValue 'D51' '71.04'
$ws.Range('E51').Value = '  -2.80%  '
